# Fruta / hortaliza, semanal
# Inserts two new weekly data rows (Castle Brite, bins 500 kilos) right after
# the first data row, pushing the existing rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows starting at row 3 (shifts rows 3:49 -> 5:51)
$ws.Rows.Item(3).Resize(2).Insert()

# --- New row 3 ---
$ws.Cells.Item(3, 1).Value = 6
$ws.Cells.Item(3, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(3, 3).Value = "Metropolitana"
$ws.Cells.Item(3, 4).Value = 44530
$ws.Cells.Item(3, 5).Value = 13
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100103
$ws.Cells.Item(3, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(3, 9).Value = 100103003
$ws.Cells.Item(3, 10).Value = "Damasco"
$ws.Cells.Item(3, 11).Value = "Castle Brite"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 14
$ws.Cells.Item(3, 14).Value = 500000
$ws.Cells.Item(3, 15).Value = 550000
$ws.Cells.Item(3, 16).Value = 525000
$ws.Cells.Item(3, 17).Value = "$/bins (500 kilos)"
$ws.Cells.Item(3, 18).Value = "Región Metropolitana"
$ws.Cells.Item(3, 19).Value = 1050
$ws.Cells.Item(3, 20).Value = 500

# --- New row 4 ---
$ws.Cells.Item(4, 1).Value = 6
$ws.Cells.Item(4, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(4, 3).Value = "Metropolitana"
$ws.Cells.Item(4, 4).Value = 44530
$ws.Cells.Item(4, 5).Value = 13
$ws.Cells.Item(4, 6).Value = "Fruta"
$ws.Cells.Item(4, 7).Value = 100103
$ws.Cells.Item(4, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(4, 9).Value = 100103003
$ws.Cells.Item(4, 10).Value = "Damasco"
$ws.Cells.Item(4, 11).Value = "Castle Brite"
$ws.Cells.Item(4, 12).Value = "Segunda"
$ws.Cells.Item(4, 13).Value = 12
$ws.Cells.Item(4, 14).Value = 400000
$ws.Cells.Item(4, 15).Value = 400000
$ws.Cells.Item(4, 16).Value = 400000
$ws.Cells.Item(4, 17).Value = "$/bins (500 kilos)"
$ws.Cells.Item(4, 18).Value = "Región Metropolitana"
$ws.Cells.Item(4, 19).Value = 800
$ws.Cells.Item(4, 20).Value = 500

# Match the date number format used by the other rows in column D (yyyy-mm-dd-style date style index)
$ws.Cells.Item(3, 4).NumberFormat = $ws.Cells.Item(5, 4).NumberFormat
$ws.Cells.Item(4, 4).NumberFormat = $ws.Cells.Item(5, 4).NumberFormat
